$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet (tab) name to reflect new "through" date
$ws.Name = "Through 2021-11-03"

# Update the November row label text
$ws.Range("A12").Value = "November (through 11-03)"

# Update November row (row 12) values for years 2016-2021 (columns C-H)
$ws.Range("C12").Value = 6
$ws.Range("D12").Value = 10
$ws.Range("E12").Value = 14
$ws.Range("F12").Value = 4
$ws.Range("G12").Value = 21
$ws.Range("H12").Value = 18

# Update Total row (row 13) values for years 2016-2021 (columns C-H)
$ws.Range("C13").Value = 492
$ws.Range("D13").Value = 720
$ws.Range("E13").Value = 629
$ws.Range("F13").Value = 486
$ws.Range("G13").Value = 1078
$ws.Range("H13").Value = 1462
